$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows right below the header (row 1) to make room for the
# newly discovered "Unnamed" columns; this pushes the existing prop_* / uf /
# partido / orgao rows down by 3.
$ws.Range("A2:A4").EntireRow.Insert()

# New rows (now at positions 2-4) for the newly found columns.
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "Unnamed: 0"

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "Unnamed: 0.1"

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "Unnamed: 0.1.1"

# Copy the formatting used by the other data rows onto the newly inserted ones.
$ws.Range("A5").Copy()
$ws.Range("A2:A4").PasteSpecial(-4122)

# Renumber the column-index values (A) for the rows that got shifted down,
# continuing the sequence started above.
$ws.Range("A5").Value = 3
$ws.Range("A6").Value = 4
$ws.Range("A7").Value = 5
$ws.Range("A8").Value = 6
$ws.Range("A9").Value = 7
$ws.Range("A10").Value = 8

# Append 3 more rows at the bottom of the table for the columns that were
# pushed past the previous end of the list.
$ws.Range("A10").Copy()
$ws.Range("A11:A13").PasteSpecial(-4122)

$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "uf"

$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "partido"

$ws.Range("A13").Value = 11
$ws.Range("B13").Value = "orgao"
